$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "protocol" sheet: rename task placeholders to the new BIDS task-naming
#    convention, replace the helper formula in G7 with literal guidance text,
#    and update the "copy exactly" instruction text.
# ---------------------------------------------------------------------------
$wsProtocol = $wb.Worksheets.Item("protocol")

# B5:B8 held the raw task name "emp" - switch to the "task-emp_bold" naming.
$wsProtocol.Cells.Item(5, 2).Value = "task-emp_bold"
$wsProtocol.Cells.Item(6, 2).Value = "task-emp_bold"
$wsProtocol.Cells.Item(7, 2).Value = "task-emp_bold"
$wsProtocol.Cells.Item(8, 2).Value = "task-emp_bold"

# G7 used to be "=tasklist!B2" (showing "emp") - replace with literal
# instructional text showing the naming pattern instead of a live formula.
$wsProtocol.Cells.Item(7, 7).Value = "task-TASKNAME_bold"

# G2 instruction text gains the TASKNAME caveat.
$wsProtocol.Cells.Item(2, 7).Value = "copy exactly and replace TASKNAME if needed"

# Wrap the long instructional text in row 2 and grow the row to fit it.
$wsProtocol.Range("A2:H2").WrapText = $true
$wsProtocol.Rows.Item(2).RowHeight = 48

# Column G needs to be wider to show the longer instructional text.
$wsProtocol.Columns.Item(7).ColumnWidth = 19.166666666666668

# ---------------------------------------------------------------------------
# 2. "participants" sheet: a folder was missing from the subject list, so a
#    new participant row is inserted (catches the missing-folder case).
# ---------------------------------------------------------------------------
$wsParticipants = $wb.Worksheets.Item("participants")
$wsParticipants.Rows.Item(4).Insert()
$wsParticipants.Cells.Item(4, 1).Value = "20170608_152333"
$wsParticipants.Cells.Item(4, 2).Value = 105
$wsParticipants.Cells.Item(4, 3).Value = "F"
$wsParticipants.Cells.Item(4, 4).Value = 18

# ---------------------------------------------------------------------------
# 3. Update selections / active sheet to match the state left behind by the
#    edit session.
# ---------------------------------------------------------------------------
$wsProtocol.Activate()
$wsProtocol.Range("G13").Select()

$wsTasklist = $wb.Worksheets.Item("tasklist")
$wsTasklist.Activate()
$wsTasklist.Range("A5").Select()

$wsParticipants.Activate()
$wsParticipants.Range("A10").Select()
